$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: updated values
$ws.Range("D2").Value = 2471
$ws.Range("E2").Value = 330
$ws.Range("F2").Value = 330
$ws.Range("G2").Value = 261
$ws.Range("H2").Value = 191
$ws.Range("I2").Value = 188
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 3518
$ws.Range("L2").Value = 2013
$ws.Range("M2").Value = 1505
$ws.Range("N2").Value = 1473
$ws.Range("O2").Value = 32
$ws.Range("P2").Value = 53
$ws.Range("Q2").Value = 226
$ws.Range("R2").Value = -522
$ws.Range("S2").Value = 338
$ws.Range("T2").Value = 565
$ws.Range("U2").Value = -339
$ws.Range("V2").Value = 1620
$ws.Range("W2").Value = 13.34
$ws.Range("X2").Value = 7.73
$ws.Range("Y2").Value = 13.62
$ws.Range("Z2").Value = 5.94
$ws.Range("AA2").Value = 133.79
$ws.Range("AB2").Value = 2590.93
$ws.Range("AC2").Value = 746
$ws.Range("AE2").Value = 5844
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 25200943
# Row 2: removed cells
$ws.Range("AD2").ClearContents()
$ws.Range("AH2").ClearContents()

# Row 3: updated values
$ws.Range("D3").Value = 2278
$ws.Range("E3").Value = 121
$ws.Range("F3").Value = 121
$ws.Range("G3").Value = 15
$ws.Range("H3").Value = 14
$ws.Range("I3").Value = 16
$ws.Range("J3").Value = -3
$ws.Range("K3").Value = 3211
$ws.Range("L3").Value = 1914
$ws.Range("M3").Value = 1297
$ws.Range("N3").Value = 1289
$ws.Range("O3").Value = 28
$ws.Range("P3").Value = 59
$ws.Range("Q3").Value = 132
$ws.Range("R3").Value = -64
$ws.Range("S3").Value = -134
$ws.Range("T3").Value = 256
$ws.Range("U3").Value = -124
$ws.Range("V3").Value = 1253
$ws.Range("W3").Value = 5.31
$ws.Range("X3").Value = 0.6
$ws.Range("Y3").Value = 1.18
$ws.Range("Z3").Value = 0.41
$ws.Range("AA3").Value = 147.54
$ws.Range("AB3").Value = 2033.49
$ws.Range("AC3").Value = 59
$ws.Range("AD3").Value = 66.84
$ws.Range("AE3").Value = 4611
$ws.Range("AF3").Value = 0.86
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 28050716

# Row 4: updated values
$ws.Range("D4").Value = 2043
$ws.Range("E4").Value = 55
$ws.Range("F4").Value = 55
$ws.Range("G4").Value = -114
$ws.Range("H4").Value = -186
$ws.Range("I4").Value = -187
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2924
$ws.Range("L4").Value = 1857
$ws.Range("M4").Value = 1067
$ws.Range("N4").Value = 1058
$ws.Range("O4").Value = 9
$ws.Range("P4").Value = 59
$ws.Range("Q4").Value = 126
$ws.Range("R4").Value = -16
$ws.Range("S4").Value = -67
$ws.Range("T4").Value = 124
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 1226
$ws.Range("W4").Value = 2.68
$ws.Range("X4").Value = -9.119999999999999
$ws.Range("Y4").Value = -15.94
$ws.Range("Z4").Value = -6.07
$ws.Range("AA4").Value = 174.09
$ws.Range("AB4").Value = 1725.52
$ws.Range("AC4").Value = -667
$ws.Range("AD4").Value = -4.03
$ws.Range("AE4").Value = 3785
$ws.Range("AF4").Value = 0.71
$ws.Range("AG4").Value = 42
$ws.Range("AH4").Value = 1.58
$ws.Range("AI4").Value = -6.3
$ws.Range("AJ4").Value = 28050716

# Row 5: updated values
$ws.Range("D5").Value = 1654
$ws.Range("E5").Value = -42
$ws.Range("F5").Value = -42
$ws.Range("G5").Value = -344
$ws.Range("H5").Value = -337
$ws.Range("I5").Value = -338
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2424
$ws.Range("L5").Value = 1742
$ws.Range("M5").Value = 683
$ws.Range("N5").Value = 673
$ws.Range("O5").Value = 9
$ws.Range("P5").Value = 140
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 15
$ws.Range("S5").Value = -89
$ws.Range("T5").Value = 87
$ws.Range("U5").Value = -87
$ws.Range("V5").Value = 1215
$ws.Range("W5").Value = -2.55
$ws.Range("X5").Value = -20.4
$ws.Range("Y5").Value = -39.04
$ws.Range("Z5").Value = -12.62
$ws.Range("AA5").Value = 255.16
$ws.Range("AB5").Value = 423.5
$ws.Range("AC5").Value = -1205
$ws.Range("AD5").Value = -1.67
$ws.Range("AE5").Value = 2430
$ws.Range("AF5").Value = 0.83
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 28050716
# Row 5: removed cells
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()

# Row 6: updated values
$ws.Range("D6").Value = 1208
$ws.Range("E6").Value = -270
$ws.Range("F6").Value = -270
$ws.Range("G6").Value = -436
$ws.Range("H6").Value = -440
$ws.Range("I6").Value = -440
$ws.Range("K6").Value = 1721
$ws.Range("L6").Value = 1440
$ws.Range("M6").Value = 281
$ws.Range("N6").Value = 272
$ws.Range("P6").Value = 147
$ws.Range("Q6").Value = 151
$ws.Range("R6").Value = 77
$ws.Range("S6").Value = -226
$ws.Range("T6").Value = 39
$ws.Range("U6").Value = 113
$ws.Range("V6").Value = 942
$ws.Range("W6").Value = -22.35
$ws.Range("X6").Value = -36.39
$ws.Range("Y6").Value = -93.06999999999999
$ws.Range("Z6").Value = -21.21
$ws.Range("AA6").Value = 512.13
$ws.Range("AB6").Value = 136.38
$ws.Range("AC6").Value = -1506
$ws.Range("AD6").Value = -1.99
$ws.Range("AE6").Value = 935
$ws.Range("AF6").Value = 3.2
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 29402067
# Row 6: removed cells
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7: removed cells
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8: removed cells
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9: removed cells
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
